$d = $word.ActiveDocument

# --- Text edits -------------------------------------------------------
# Replace the company-name placeholder with the real team name.
$d.Content.Find.Execute("<Insert_Company_Name>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Jayhawks", 2) | Out-Null

# "Present Members: ..." was split across two runs ("...Victor, " + "Ellia").
# Re-asserting the full text through Find/Replace collapses it back into a
# single run, matching the committed version.
$d.Content.Find.Execute("Present Members: Timo, Allie, Riley, Victor, Ellia",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "Present Members: Timo, Allie, Riley, Victor, Ellia", 2) | Out-Null

# --- Style edits --------------------------------------------------------
# Both the document's default paragraph style ("Normal") and the LibreOffice
# "LO-normal" style (internally "Normal1") gained an explicit
# suppressAutoHyphens setting. Word exposes this through
# Style.ParagraphFormat.Hyphenation (True == hyphenation allowed ==
# suppressAutoHyphens absent/false; False == suppressAutoHyphens true).
$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.Hyphenation = $false

$loNormal = $d.Styles.Item("LO-normal")
$loNormal.ParagraphFormat.Hyphenation = $false
